# Fruta / hortaliza, semanal
# Insert this week's new price observations (2022-02-16) for
# "Pepino ensalada" at row 31/32, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 31 (shifts existing rows 31.. down by 2)
$ws.Rows.Item(31).EntireRow.Insert()
$ws.Rows.Item(31).EntireRow.Insert()

# Row 31: Primera quality, new weekly observation
$ws.Cells.Item(31, 1).Value = 2
$ws.Cells.Item(31, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = [DateTime]"2022-02-16"
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112043
$ws.Cells.Item(31, 7).Value = "Pepino ensalada"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 500
$ws.Cells.Item(31, 11).Value = 11000
$ws.Cells.Item(31, 12).Value = 12000
$ws.Cells.Item(31, 13).Value = 11500
$ws.Cells.Item(31, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 164
$ws.Cells.Item(31, 17).Value = 70
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Row 32: Segunda quality, new weekly observation
$ws.Cells.Item(32, 1).Value = 2
$ws.Cells.Item(32, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = [DateTime]"2022-02-16"
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100112043
$ws.Cells.Item(32, 7).Value = "Pepino ensalada"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Segunda"
$ws.Cells.Item(32, 10).Value = 400
$ws.Cells.Item(32, 11).Value = 9000
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 9500
$ws.Cells.Item(32, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 95
$ws.Cells.Item(32, 17).Value = 100
$ws.Cells.Item(32, 18).Value = "Hortaliza"
